# Updated cryptos list (GitHub Actions refresh): prices / 1h volumes, and
# the ARBITRUM / VeChain rows 38-39 swapped places with refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.023.12'
$ws.Range("E2").Value = '  +2.71%  '

$ws.Range("D3").Value = '1.651.37'
$ws.Range("E3").Value = '  +3.44%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.62%  '

$ws.Range("E6").Value = '  +1.39%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.249'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.74%  '

$ws.Range("E9").Value = '  +1.61%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.86'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0866'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.20%  '

$ws.Range("D12").Value = '1.886.15'

$ws.Range("D13").Value = '1.636.09'
$ws.Range("E13").Value = '  +2.64%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.07'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.15%  '

$ws.Range("E15").Value = '  +2.81%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.31'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.88%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '239.94'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.03%  '

$ws.Range("D18").Value = '27.004.99'
$ws.Range("E18").Value = '  +2.67%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.67%  '

$ws.Range("E20").Value = '  +1.22%  '

$ws.Range("E21").Value = '  +0.01%  '

$ws.Range("E22").Value = '  +4.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.07%  '

$ws.Range("E24").Value = '  +3.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.85'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.56%  '

$ws.Range("E26").Value = '  -0.11%  '

$ws.Range("E27").Value = '  +2.02%  '

$ws.Range("E28").Value = '  +1.86%  '

$ws.Range("E29").Value = '  +2.89%  '

$ws.Range("E30").Value = '  +0.57%  '

$ws.Range("E31").Value = '  +1.76%  '

$ws.Range("E32").Value = '  +3.31%  '

$ws.Range("D33").Value = '1.523.08'
$ws.Range("E33").Value = '  +0.51%  '

$ws.Range("E34").Value = '  +5.12%  '

$ws.Range("E35").Value = '  +8.62%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.42'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.08%  '

$ws.Range("E37").Value = '  +2.47%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0169'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.99%  '

$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.885'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.43%  '

$ws.Range("E40").Value = '  +2.78%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.26'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.77'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.75%  '

$ws.Range("D44").Value = '1.792.02'
$ws.Range("E44").Value = '  +3.31%  '

$ws.Range("E45").Value = '  +2.10%  '

$ws.Range("E46").Value = '  -2.18%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '89.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.32%  '

$ws.Range("E48").Value = '  +0.89%  '

$ws.Range("E49").Value = '  +2.61%  '

$ws.Range("E50").Value = '  +1.34%  '

$ws.Range("E51").Value = '  +2.01%  '

